$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): Wins, Losses, Ties in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (AC1) to new header cells so they match
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2-45: Wins=86, Losses=76, Ties=0
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD
    $ws.Cells.Item($r, 31).Value = 76   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
